# Booking.cshtml und Booking im BookingController erstellt
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Row 43: "Booking Methode im BookingController erstellen" is now done,
# with D43 becoming a real date (2019-04-02) instead of the placeholder
# text "bis 28.03.2019". Copy D41's date formatting (reuses the existing
# date style) then set the new value.
$ws.Range("D41").Copy()
$ws.Range("D43").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D43").Value = "04/02/2019"
$ws.Range("C43").Value = "done"

# Row 44: "Booking View erstellen" is now done too.
$ws.Range("D41").Copy()
$ws.Range("D44").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D44").Value = "04/02/2019"
$ws.Range("C44").Value = "done"

$excel.CutCopyMode = $false

# Insert two new rows right after row 44 for the BookingConfirmation tasks,
# pushing "email Bestätigung einbauen" from row 45 down to row 47.
$ws.Range("A45:A46").EntireRow.Insert()
# The insert carries D44's date formatting down into the new rows as blank
# styled cells; drop that so 45/46 only carry the B-column task text.
$ws.Range("C45:D46").Clear()

$ws.Range("B45").Value = "BookingConfirmation Methode im BookingController erstellen"
$ws.Range("B46").Value = "BookingConfirmation View erstellen"

# Update the active selection to match the saved workbook state.
$ws.Range("C45").Select()
